# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 09:53"

# Row 6 - India
$ws.Range("B6").Value = 1807222
$ws.Range("C6").Value = 2520
$ws.Range("E6").Value = 580657

# Row 7 - Rusia
$ws.Range("B7").Value = 856264
$ws.Range("C7").Value = 5394
$ws.Range("D7").Value = 653593
$ws.Range("E7").Value = 188464
$ws.Range("G7").Value = 79
$ws.Range("H7").Value = 14207

# Row 46 - Rumania
$ws.Range("B46").Value = 53051
$ws.Range("C46").Value = 226
$ws.Range("E46").Value = 6098

# Row 65 - Argelia
$ws.Range("D65").Value = 17942
$ws.Range("E65").Value = 6628
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 792

# Row 87 - Sudan
$ws.Range("B87").Value = 9133
$ws.Range("C87").Value = 18
$ws.Range("D87").Value = 7424
$ws.Range("E87").Value = 1494

# Row 94 - Guayana Francesa
$ws.Range("B94").Value = 7364
$ws.Range("C94").Value = 47
$ws.Range("D94").Value = 6505
$ws.Range("E94").Value = 813

# Row 124 - Suazilandia
$ws.Range("B124").Value = 2354
$ws.Range("C124").Value = 10
$ws.Range("D124").Value = 1746
$ws.Range("E124").Value = 579

# Row 127 - Mali
$ws.Range("B127").Value = 2120
$ws.Range("C127").Value = 10
$ws.Range("D127").Value = 1645
$ws.Range("E127").Value = 395

# Row 128 - Sudan del Sur
$ws.Range("B128").Value = 2080
$ws.Range("C128").Value = 1
$ws.Range("D128").Value = 1935

# Row 140 - Benin
$ws.Range("B140").Value = 1246
$ws.Range("C140").Value = 3
$ws.Range("E140").Value = 162
